$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 859.61536
$ws.Range("I18").Value = 336.44446
$ws.Range("J18").Value = 2036.75
$ws.Range("K18").Value = 336.44446
$ws.Range("L18").Value = 2036.75
$ws.Range("M18").Value = -52.44445999999999
$ws.Range("N18").Value = -2604.75

$ws.Range("H55").Value = 92113
$ws.Range("I55").Value = 339.5
$ws.Range("K55").Value = 339.5
$ws.Range("M55").Value = -125.5

$ws.Range("H112").Value = 3010.3333
$ws.Range("J112").Value = 3150.3872
$ws.Range("L112").Value = 9451.161599999999
$ws.Range("N112").Value = -11667.1616

$ws.Range("H131").Value = 2368.516
$ws.Range("I131").Value = 1439.4166
$ws.Range("J131").Value = 5554
$ws.Range("K131").Value = 4318.2498
$ws.Range("L131").Value = 16662
$ws.Range("M131").Value = 721.7502000000004
$ws.Range("N131").Value = -26742

$ws.Range("H138").Value = 5261.18
$ws.Range("I138").Value = 2463.2
$ws.Range("J138").Value = 5960.675
$ws.Range("K138").Value = 7389.599999999999
$ws.Range("L138").Value = 17882.025
$ws.Range("M138").Value = -2249.599999999999
$ws.Range("N138").Value = -28162.025

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 50002080
$ws.Range("J20").Value = 2320
$ws.Range("L20").Value = 2320
$ws.Range("N20").Value = -2814

$ws.Range("H105").Value = 2819.3242
$ws.Range("I105").Value = 2903.8333
$ws.Range("J105").Value = 2778.76
$ws.Range("K105").Value = 2903.8333
$ws.Range("L105").Value = 2778.76
$ws.Range("M105").Value = -1156.8333
$ws.Range("N105").Value = -6272.76

$ws.Range("H134").Value = 4903.278
$ws.Range("I134").Value = 3195.2727
$ws.Range("K134").Value = 9585.8181
$ws.Range("M134").Value = -7050.8181

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 560013.8
$ws.Range("I58").Value = 1001803
$ws.Range("K58").Value = 1001803
$ws.Range("M58").Value = -1001600

$ws.Range("H62").Value = 3218.8235
$ws.Range("J62").Value = 3119.5557
$ws.Range("L62").Value = 3119.5557
$ws.Range("N62").Value = -4367.5557

$ws.Range("H65").Value = 3218.8235
$ws.Range("J65").Value = 3119.5557
$ws.Range("L65").Value = 15597.7785
$ws.Range("N65").Value = -21837.7785

$ws.Range("H136").Value = 560013.8
$ws.Range("I136").Value = 1001803
$ws.Range("K136").Value = 3005409
$ws.Range("M136").Value = -3002859

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3731.5
$ws.Range("I3").Value = 1642
$ws.Range("K3").Value = 4926
$ws.Range("M3").Value = -4814

$ws.Range("H7").Value = 1102
$ws.Range("J7").Value = 1163.3334
$ws.Range("L7").Value = 3490.0002
$ws.Range("N7").Value = -3714.0002

$ws.Range("H70").Value = 1666.3334
$ws.Range("I70").Value = 1999.5
$ws.Range("K70").Value = 5998.5
$ws.Range("M70").Value = -5683.5

$ws.Range("H73").Value = 1666.3334
$ws.Range("I73").Value = 1999.5
$ws.Range("K73").Value = 5998.5
$ws.Range("M73").Value = -4906.5

$ws.Range("H92").Value = 1030
$ws.Range("J92").Value = 1390.8334
$ws.Range("L92").Value = 4172.5002
$ws.Range("N92").Value = -6668.5002

$ws.Range("H117").Value = 437.66666
$ws.Range("I117").Value = 437.66666
$ws.Range("K117").Value = 1312.99998
$ws.Range("M117").Value = 2129.00002

$ws.Range("H127").Value = 15030
$ws.Range("I127").Value = 10030
$ws.Range("J127").Value = 20030
$ws.Range("K127").Value = 30090
$ws.Range("L127").Value = 60090
$ws.Range("M127").Value = -25130
$ws.Range("N127").Value = -70010

$ws.Range("H131").Value = 5396.6875
$ws.Range("I131").Value = 1690.2858
$ws.Range("K131").Value = 5070.857400000001
$ws.Range("M131").Value = -30.85740000000078

$ws.Range("H132").Value = 3789.818
$ws.Range("J132").Value = 3818.8
$ws.Range("L132").Value = 34369.2
$ws.Range("N132").Value = -39429.2

$ws.Range("H134").Value = 10272.429
$ws.Range("I134").Value = 9484.5
$ws.Range("K134").Value = 28453.5
$ws.Range("M134").Value = -23383.5

$ws.Range("H138").Value = 2479.75
$ws.Range("I138").Value = 973.3333
$ws.Range("J138").Value = 6999
$ws.Range("K138").Value = 2919.9999
$ws.Range("L138").Value = 20997
$ws.Range("M138").Value = 2220.0001
$ws.Range("N138").Value = -31277

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 18359
$ws.Range("J63").Value = 18359
$ws.Range("L63").Value = 18359
$ws.Range("N63").Value = -19731

$ws.Range("H66").Value = 18359
$ws.Range("J66").Value = 18359
$ws.Range("L66").Value = 55077
$ws.Range("N66").Value = -61941

$ws.Range("H80").Value = 2503002
$ws.Range("I80").Value = 3335668
$ws.Range("K80").Value = 3335668
$ws.Range("M80").Value = -3334670

$ws.Range("H83").Value = 2503002
$ws.Range("I83").Value = 3335668
$ws.Range("K83").Value = 16678340
$ws.Range("M83").Value = -16673348

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 912.5
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 912.5
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws.Range("H46").Value = 3253.2354
$ws.Range("I46").Value = 3009.2727
$ws.Range("J46").Value = 3700.5
$ws.Range("K46").Value = 3009.2727
$ws.Range("L46").Value = 3700.5
$ws.Range("M46").Value = -2821.2727
$ws.Range("N46").Value = -4076.5

$ws.Range("H55").Value = 989.0909
$ws.Range("J55").Value = 2033
$ws.Range("L55").Value = 2033
$ws.Range("N55").Value = -2379

$ws.Range("H68").Value = 2045.1818
$ws.Range("I68").Value = 1749.7
$ws.Range("K68").Value = 1749.7
$ws.Range("M68").Value = -1000.7

$ws.Range("H71").Value = 2045.1818
$ws.Range("I71").Value = 1749.7
$ws.Range("K71").Value = 8748.5
$ws.Range("M71").Value = -5004.5

$ws.Range("H82").Value = 3981.1052
$ws.Range("I82").Value = 3458.4546
$ws.Range("J82").Value = 4699.75
$ws.Range("K82").Value = 3458.4546
$ws.Range("L82").Value = 4699.75
$ws.Range("M82").Value = -3097.4546
$ws.Range("N82").Value = -5421.75

$ws.Range("H85").Value = 3981.1052
$ws.Range("I85").Value = 3458.4546
$ws.Range("J85").Value = 4699.75
$ws.Range("K85").Value = 3458.4546
$ws.Range("L85").Value = 4699.75
$ws.Range("M85").Value = -2210.4546
$ws.Range("N85").Value = -7195.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 731.1429000000001
$ws.Range("I2").Value = 683
$ws.Range("J2").Value = 851.5
$ws.Range("K2").Value = 683
$ws.Range("L2").Value = 851.5
$ws.Range("M2").Value = -571
$ws.Range("N2").Value = -1075.5

$ws.Range("H54").Value = 33019
$ws.Range("J54").Value = 33019
$ws.Range("L54").Value = 33019
$ws.Range("N54").Value = -34059

$ws.Range("H62").Value = 84846.5
$ws.Range("I62").Value = 130550.25
$ws.Range("K62").Value = 130550.25
$ws.Range("M62").Value = -129926.25

$ws.Range("H65").Value = 84846.5
$ws.Range("I65").Value = 130550.25
$ws.Range("K65").Value = 652751.25
$ws.Range("M65").Value = -649631.25

$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("N90").ClearContents()

$ws.Range("H96").Value = 102935
$ws.Range("I96").Value = 112710.78
$ws.Range("J96").Value = 14953
$ws.Range("K96").Value = 112710.78
$ws.Range("L96").Value = 14953
$ws.Range("M96").Value = -111337.78
$ws.Range("N96").Value = -17699

$ws.Range("H132").Value = 3057.0466
$ws.Range("I132").Value = 2184.3713
$ws.Range("J132").Value = 6875
$ws.Range("K132").Value = 6553.113899999999
$ws.Range("L132").Value = 20625
$ws.Range("M132").Value = -4023.113899999999
$ws.Range("N132").Value = -25685
